# 22/03 - LMA: New User Story not yet in scope. Minor Fixes
#
# Populates the "D03NonQuotableProducts" sheet (rows B3:B66) with the list of
# non-quotable products for the new user story, and makes that sheet the
# active tab/selection (it was previously on "PhoneLine").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("D03NonQuotableProducts")

# Switch focus to this sheet (this is the tab that ends up active/selected
# in the saved workbook, replacing the previously active "PhoneLine" tab).
$ws.Activate()

$ws.Range("B3").Value = "ONE Project"
$ws.Range("B4").Value = "Explore"
$ws.Range("B5").Value = "E-Line"
$ws.Range("B6").Value = "Explore BiLAN Teleworking"
$ws.Range("B7").Value = "Explore International with voice"
$ws.Range("B8").Value = "Explore Mobile Worker"
$ws.Range("B9").Value = "Explore Mono with voice"
$ws.Range("B10").Value = "Explore Mono without voice"
$ws.Range("B11").Value = "SDWAN"
$ws.Range("B12").Value = "SDWAN International"
$ws.Range("B13").Value = "IP Pack on BiLAN/Explore"
$ws.Range("B14").Value = "Microsoft Office 365 MS Teams"
$ws.Range("B15").Value = "Temporary rental PABX"
$ws.Range("B16").Value = "Call Connect"
$ws.Range("B17").Value = "PRA over IAD30"
$ws.Range("B18").Value = "Directory Number (DDI/ISDN)"
$ws.Range("B19").Value = "PRA"
$ws.Range("B20").Value = "Bizz Call Connect"
$ws.Range("B21").Value = "Multiline"
$ws.Range("B22").Value = "Cloud Mail Security"
$ws.Range("B23").Value = "DNS"
$ws.Range("B24").Value = "Personal Cloud"
$ws.Range("B25").Value = "Tariff Plan"
$ws.Range("B26").Value = "Cloud Exchange"
$ws.Range("B27").Value = "Cloud Mail Security"
$ws.Range("B28").Value = "Conversational Chatbot"
$ws.Range("B29").Value = "EM+S (Microsoft)"
$ws.Range("B30").Value = "Hosting (Shared)"
$ws.Range("B31").Value = "Housing (not BiLAN)"
$ws.Range("B32").Value = "Interact"
$ws.Range("B33").Value = "Internet of Things"
$ws.Range("B34").Value = "Microsoft Office 365"
$ws.Range("B35").Value = "Proximus Azure Services"
$ws.Range("B36").Value = "Proximus DocDrop"
$ws.Range("B37").Value = "Secure Mail"
$ws.Range("B38").Value = "SMS Solution"
$ws.Range("B39").Value = "SMS Solutions Pack"
$ws.Range("B40").Value = "Voice Assist"
$ws.Range("B41").Value = "Workspace_One (Airwatch)"
$ws.Range("B42").Value = "Cloud vContainer"
$ws.Range("B43").Value = "Temporary xDSL Fast Internet"
$ws.Range("B44").Value = "Marketing Number"
$ws.Range("B45").Value = "Marketing Number International"
$ws.Range("B46").Value = "Marketing Number Mobile"
$ws.Range("B47").Value = "VMS"
$ws.Range("B48").Value = "Mass Market"
$ws.Range("B49").Value = "Pack (with mobile)"
$ws.Range("B50").Value = "Pack (without mobile)"
$ws.Range("B51").Value = "Mobile Prepaid"
$ws.Range("B52").Value = "Mobile BONE"
$ws.Range("B53").Value = "Elevator line"
$ws.Range("B54").Value = "Temporary ISDN"
$ws.Range("B55").Value = "Engage"
$ws.Range("B56").Value = "Leased Line"
$ws.Range("B57").Value = "CPE"
$ws.Range("B58").Value = "Split Plan-IFE/PFE"
$ws.Range("B59").Value = "Joint Offer"
$ws.Range("B60").Value = "Proximus TV"
$ws.Range("B61").Value = "Proximus TV App (TV Everywhere)"
$ws.Range("B62").Value = "Split Plan/IFE Tool"
$ws.Range("B63").Value = "M2M KORE"
$ws.Range("B64").Value = "Fixed IP address"
$ws.Range("B65").Value = "Tariff Plan"
$ws.Range("B66").Value = "IP Pack on xDSL Fast Internet"

# Final cursor position on this sheet, matching the saved selection.
$ws.Range("B19").Select() | Out-Null
